$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 10 to make room for the Gaussian-Quadrature
# entry (moved up from the bottom of the table) and the 3 new Spiral-scheme
# entries that were added.
$ws.Rows("10:12").Insert()

# The 3 freshly inserted rows need column A formatted like the rest of the
# "#" column (bold, centered, thin box border) to match the other data rows.
$ws.Range("A10:A12").Font.Bold = $true
$ws.Range("A10:A12").HorizontalAlignment = -4108
$ws.Range("A10:A12").VerticalAlignment = -4160
$ws.Range("A10:A12").Borders.LineStyle = 1

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.225040886154097
$ws.Range("D10").Value = 2.34493493238325
$ws.Range("E10").Value = 0.6296575682396935
$ws.Range("F10").Value = 1.225040886154097
$ws.Range("G10").Value = 1.332156129456757
$ws.Range("H10").Value = 0.4592399385408076
$ws.Range("I10").Value = 0.7660167035182626
$ws.Range("J10").Value = 2.34493493238325
$ws.Range("K10").Value = 1.487296250311472
$ws.Range("L10").Value = 1.356168568232785
$ws.Range("M10").Value = 1.126174359715478

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.7971362181471839
$ws.Range("D11").Value = 0.9376980776736024
$ws.Range("E11").Value = 1.178874050533731
$ws.Range("F11").Value = 0.7971362181471839
$ws.Range("G11").Value = 0.7221209008683771
$ws.Range("H11").Value = 1.86278677253059
$ws.Range("I11").Value = 1.042002421923074
$ws.Range("J11").Value = 0.9376980776736024
$ws.Range("K11").Value = 1.058286064103666
$ws.Range("L11").Value = 0.9277111411254252
$ws.Range("M11").Value = 1.09010307361276

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.7905430295222121
$ws.Range("D12").Value = 0.9385155587952247
$ws.Range("E12").Value = 1.18069548918641
$ws.Range("F12").Value = 0.7905430295222121
$ws.Range("G12").Value = 0.7232876998929499
$ws.Range("H12").Value = 1.865748641074372
$ws.Range("I12").Value = 1.04175520716646
$ws.Range("J12").Value = 0.9385155587952247
$ws.Range("K12").Value = 1.059605523990817
$ws.Range("L12").Value = 0.9250742767565147
$ws.Range("M12").Value = 1.090090937606272

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.7960958519844883
$ws.Range("D13").Value = 0.9362009970108627
$ws.Range("E13").Value = 1.17910867606748
$ws.Range("F13").Value = 0.7960958519844883
$ws.Range("G13").Value = 0.722567982635019
$ws.Range("H13").Value = 1.862021717421377
$ws.Range("I13").Value = 1.041868783256412
$ws.Range("J13").Value = 0.9362009970108627
$ws.Range("K13").Value = 1.057654836539171
$ws.Range("L13").Value = 0.9268753442618298
$ws.Range("M13").Value = 1.08964400139594

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.1309559999999998
$ws.Range("D14").Value = 4.350464000000007
$ws.Range("E14").Value = 0.6593600000000008
$ws.Range("F14").Value = 0.1309559999999998
$ws.Range("G14").Value = 1.81494
$ws.Range("H14").Value = 0.8801960000000019
$ws.Range("I14").Value = 0.4860880000000009
$ws.Range("J14").Value = 4.350464000000007
$ws.Range("K14").Value = 2.504912000000004
$ws.Range("L14").Value = 1.317934000000001
$ws.Range("M14").Value = 1.387000666666668

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.01
$ws.Range("D15").Value = 6.875437500000012
$ws.Range("E15").Value = 0.15
$ws.Range("F15").Value = 0.01
$ws.Range("G15").Value = 2.526812499999997
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.11
$ws.Range("J15").Value = 6.875437500000012
$ws.Range("K15").Value = 3.512718750000006
$ws.Range("L15").Value = 1.761359375000003
$ws.Range("M15").Value = 1.612041666666668

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.4217035024383986
$ws.Range("D16").Value = 4.380577564671997
$ws.Range("E16").Value = 0.4936157308928026
$ws.Range("F16").Value = 0.4217035024383986
$ws.Range("G16").Value = 1.888129161625597
$ws.Range("H16").Value = 0.4151164925952002
$ws.Range("I16").Value = 0.4786240591872011
$ws.Range("J16").Value = 4.380577564671997
$ws.Range("K16").Value = 2.4370966477824
$ws.Range("L16").Value = 1.429400075110399
$ws.Range("M16").Value = 1.346294418568533

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.996822371435004
$ws.Range("D17").Value = 0.9933785644704798
$ws.Range("E17").Value = 0.993236875896363
$ws.Range("F17").Value = 0.996822371435004
$ws.Range("G17").Value = 0.9961857860851996
$ws.Range("H17").Value = 0.98684864490651
$ws.Range("I17").Value = 0.993206737851276
$ws.Range("J17").Value = 0.9933785644704798
$ws.Range("K17").Value = 0.9933077201834214
$ws.Range("L17").Value = 0.9950650458092126
$ws.Range("M17").Value = 0.993279830107472

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.001146461615077
$ws.Range("D18").Value = 0.8707963194674984
$ws.Range("E18").Value = 1.035743388192315
$ws.Range("F18").Value = 1.001146461615077
$ws.Range("G18").Value = 0.9408042698655251
$ws.Range("H18").Value = 1.070568582324037
$ws.Range("I18").Value = 1.020470608105398
$ws.Range("J18").Value = 0.8707963194674984
$ws.Range("K18").Value = 0.9532698538299068
$ws.Range("L18").Value = 0.977208157722492
$ws.Range("M18").Value = 0.9899216049283085

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.132588323769258
$ws.Range("D19").Value = 0.729469432657212
$ws.Range("E19").Value = 1.034511192148362
$ws.Range("F19").Value = 1.132588323769258
$ws.Range("G19").Value = 0.7858181002285244
$ws.Range("H19").Value = 1.251594718786602
$ws.Range("I19").Value = 1.07169064935985
$ws.Range("J19").Value = 0.729469432657212
$ws.Range("K19").Value = 0.8819903124027868
$ws.Range("L19").Value = 1.007289318086022
$ws.Range("M19").Value = 1.000945402824968

